# Remove the Status/Remarks/Actual Output/Screenshot/Page Source results
# (columns L:P) for rows 10-24, and Status/Remarks (columns L:M) for rows
# 25-40 on the single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10-24: clear columns L through P
$ws.Range("L10:P24").ClearContents()

# Rows 25-40: clear columns L through M
$ws.Range("L25:M40").ClearContents()
